# Apply cryptos list update (prices / volume % / a couple of coin row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.937.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.46"
$ws.Range("D3").ClearFormats()

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.91"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.40%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.36%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.13"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.11%  "

# Row 11
$ws.Range("E11").Value = "  +3.73%  "

# Row 13
$ws.Range("E13").Value = "  -3.35%  "

# Row 14
$ws.Range("E14").Value = "  -0.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.729.46"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.357.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.21%  "

# Row 17
$ws.Range("E17").Value = "  +0.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.896.27"
$ws.Range("D18").ClearFormats()

# Row 19
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.11%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.11%  "

# Row 21
$ws.Range("E21").Value = "  -0.79%  "

# Row 22
$ws.Range("E22").Value = "  +0.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.88"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("E24").Value = "  -5.05%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.41%  "

# Row 27
$ws.Range("E27").Value = "  +0.76%  "

# Row 28
$ws.Range("E28").Value = "  +0.83%  "

# Row 29
$ws.Range("E29").Value = "  +2.21%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.93"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.47%  "

# Row 31
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("E32").Value = "  +0.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.35"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0712"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.06%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.85"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.24%  "

# Row 36
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "127.58"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -23.06%  "

# Row 37
$ws.Range("E37").Value = "  +3.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.27"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.93%  "

# Row 41
$ws.Range("E41").Value = "  -0.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.19"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.41%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.931.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0278"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.27%  "

# Row 45
$ws.Range("E45").Value = "  +2.51%  "

# Row 46
$ws.Range("E46").Value = "  -0.67%  "

# Row 47
$ws.Range("E47").Value = "  -8.67%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.594.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.11%  "

# Row 49
$ws.Range("E49").Value = "  +1.73%  "

# Row 50
$ws.Range("E50").Value = "  +1.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.30%  "
